# Reproduces the TBV report content addition: 14 new paragraphs
# (tag markers, 'List Bullet' requirement headers, and indented
# detail/body paragraphs) inserted right after the existing blank
# paragraph that follows the 'TBV Tags' title.
$d = $word.ActiveDocument

# Create all of the new paragraphs first (as plain/Normal, empty
# paragraphs) before applying any per-paragraph style/indent, so that
# a 'List Bullet' style applied to one paragraph can't bleed forward
# onto paragraphs split off from it afterwards.
$cur = $d.Paragraphs.Last
$newParas = @()
for ($i = 0; $i -lt 14; $i++) {
    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Last
    $newParas += $cur
}

# 0: plain
$newParas[0].Range.InsertBefore("[PUMP:TBV:1111]")

# 1: ListBullet
$newParas[1].Style = "List Bullet"
$newParas[1].Range.InsertBefore("PUMP:HRD:3350")

# 2: indent 360
$newParas[2].LeftIndent = 18
$newParas[2].Range.InsertBefore("Details regarding the full color touchscreen. ")

# 3: indent 720
$newParas[3].LeftIndent = 36
$newParas[3].Range.InsertBefore("PUMP:HTP:1500")

# 4: indent 720
$newParas[4].LeftIndent = 36
$newParas[4].Range.InsertBefore("Test 1500 ")

# 5: ListBullet
$newParas[5].Style = "List Bullet"
$newParas[5].Range.InsertBefore("PUMP:HRD:0000")

# 6: indent 360
$newParas[6].LeftIndent = 18
$newParas[6].Range.InsertBefore("Details regarding the size and weight of the pump. ")

# 7: plain
$newParas[7].Range.InsertBefore("[PUMP:TBV:1]")

# 8: ListBullet
$newParas[8].Style = "List Bullet"
$newParas[8].Range.InsertBefore("ACE:SRS:1")

# 9: indent 360
$newParas[9].LeftIndent = 18
$newParas[9].Range.InsertBefore("The software shall provide a bolus feature which generates boluses in the range of 0.01 to 25 units, which an increment of 0.01 units.  ")

# 10: indent 720
$newParas[10].LeftIndent = 36
$newParas[10].Range.InsertBefore("PUMP:SDS:10")

# 11: indent 720
$newParas[11].LeftIndent = 36
$newParas[11].Range.InsertBefore("Here are details of how the bolus calculator works ….         ")

# 12: indent 720
$newParas[12].LeftIndent = 36
$newParas[12].Range.InsertBefore("PUMP:SVAL:100")

# 13: indent 720
$newParas[13].LeftIndent = 36
$newParas[13].Range.InsertBefore("This test validates bolus features… blah, blah, blah        ")
